$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2083333333333333
$ws.Range("C2").Value = 0.5128205128205128
$ws.Range("J2").Value = 0.03525641025641026
$ws.Range("P2").Value = 0.1217948717948718
$ws.Range("S2").Value = 0.1217948717948718
$ws.Range("B3").Value = 0.01162790697674419
$ws.Range("C3").Value = 0.03488372093023256
$ws.Range("J3").Value = 0.04651162790697674
$ws.Range("P3").Value = 0.7441860465116279
$ws.Range("S3").Value = 0.1627906976744186
$ws.Range("J4").Value = 0.0975609756097561
$ws.Range("P4").Value = 0.7317073170731707
$ws.Range("S4").Value = 0.1707317073170732
$ws.Range("B6").Value = 0.06589147286821706
$ws.Range("D6").Value = 0.007751937984496124
$ws.Range("E6").Value = 0.003875968992248062
$ws.Range("F6").Value = 0.04263565891472868
$ws.Range("J6").Value = 0.251937984496124
$ws.Range("O6").Value = 0.0310077519379845
$ws.Range("Q6").Value = 0.1511627906976744
$ws.Range("R6").Value = 0.04263565891472868
$ws.Range("S6").Value = 0.4031007751937984
$ws.Range("B7").Value = 0.07100591715976332
$ws.Range("D7").Value = 0.02366863905325444
$ws.Range("F7").Value = 0.05917159763313609
$ws.Range("J7").Value = 0.1834319526627219
$ws.Range("Q7").Value = 0.1597633136094675
$ws.Range("R7").Value = 0.0650887573964497
$ws.Range("S7").Value = 0.4378698224852071
$ws.Range("B8").Value = 0.08713692946058091
$ws.Range("D8").Value = 0.01867219917012448
$ws.Range("F8").Value = 0.06431535269709543
$ws.Range("J8").Value = 0.0975103734439834
$ws.Range("O8").Value = 0.02074688796680498
$ws.Range("Q8").Value = 0.1804979253112033
$ws.Range("R8").Value = 0.08506224066390042
$ws.Range("S8").Value = 0.4460580912863071
$ws.Range("B9").Value = 0.05161290322580645
$ws.Range("D9").Value = 0.01612903225806452
$ws.Range("E9").Value = 0.003225806451612903
$ws.Range("F9").Value = 0.07741935483870968
$ws.Range("J9").Value = 0.0935483870967742
$ws.Range("O9").Value = 0.01290322580645161
$ws.Range("Q9").Value = 0.2064516129032258
$ws.Range("R9").Value = 0.08387096774193549
$ws.Range("S9").Value = 0.4548387096774194
$ws.Range("B10").Value = 0.1266932270916335
$ws.Range("D10").Value = 0.01832669322709163
$ws.Range("F10").Value = 0.08366533864541832
$ws.Range("J10").Value = 0.1243027888446215
$ws.Range("O10").Value = 0.01593625498007968
$ws.Range("Q10").Value = 0.2143426294820717
$ws.Range("R10").Value = 0.06772908366533864
$ws.Range("S10").Value = 0.349003984063745
$ws.Range("G11").Value = 0.1264822134387352
$ws.Range("J11").Value = 0.1264822134387352
$ws.Range("K11").Value = 0.158102766798419
$ws.Range("L11").Value = 0.5731225296442688
$ws.Range("S11").Value = 0.0158102766798419
$ws.Range("G12").Value = 0.7548387096774194
$ws.Range("J12").Value = 0.1290322580645161
$ws.Range("K12").Value = 0.03870967741935484
$ws.Range("L12").Value = 0.05806451612903226
$ws.Range("S12").Value = 0.01935483870967742
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2727272727272727
$ws.Range("S13").Value = 0.06060606060606061
$ws.Range("F15").Value = 0.0273972602739726
$ws.Range("H15").Value = 0.1780821917808219
$ws.Range("I15").Value = 0.1050228310502283
$ws.Range("J15").Value = 0.3242009132420091
$ws.Range("K15").Value = 0.0410958904109589
$ws.Range("M15").Value = 0.0136986301369863
$ws.Range("N15").Value = 0.0045662100456621
$ws.Range("O15").Value = 0.0365296803652968
$ws.Range("S15").Value = 0.2694063926940639
$ws.Range("F16").Value = 0.04712041884816754
$ws.Range("H16").Value = 0.1204188481675393
$ws.Range("I16").Value = 0.09947643979057591
$ws.Range("J16").Value = 0.3507853403141361
$ws.Range("K16").Value = 0.1099476439790576
$ws.Range("M16").Value = 0.005235602094240838
$ws.Range("N16").Value = 0.01047120418848168
$ws.Range("O16").Value = 0.0418848167539267
$ws.Range("S16").Value = 0.2146596858638743
$ws.Range("F17").Value = 0.01446280991735537
$ws.Range("H17").Value = 0.1880165289256198
$ws.Range("I17").Value = 0.1260330578512397
$ws.Range("J17").Value = 0.3863636363636364
$ws.Range("K17").Value = 0.08884297520661157
$ws.Range("M17").Value = 0.01033057851239669
$ws.Range("N17").Value = 0.002066115702479339
$ws.Range("O17").Value = 0.04752066115702479
$ws.Range("S17").Value = 0.1363636363636364
$ws.Range("F18").Value = 0.02857142857142857
$ws.Range("H18").Value = 0.2057142857142857
$ws.Range("I18").Value = 0.1085714285714286
$ws.Range("J18").Value = 0.3485714285714286
$ws.Range("K18").Value = 0.12
$ws.Range("M18").Value = 0.005714285714285714
$ws.Range("O18").Value = 0.07428571428571429
$ws.Range("S18").Value = 0.1085714285714286
$ws.Range("F19").Value = 0.01401541695865452
$ws.Range("H19").Value = 0.2053258584442887
$ws.Range("I19").Value = 0.1317449194113525
$ws.Range("J19").Value = 0.3454800280308339
$ws.Range("K19").Value = 0.07848633496846531
$ws.Range("M19").Value = 0.0161177295024527
$ws.Range("N19").Value = 0.002102312543798178
$ws.Range("O19").Value = 0.06937631394533987
$ws.Range("S19").Value = 0.1373510861948143
